$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 17: new entry for the day worked on uPlot zooming
$ws.Range("A17").Value = 44354
$ws.Range("B17").Value = 7
$ws.Range("C17").Formula = "=C16+B17"
$ws.Range("D17").Value = "Implemented zooming for uPlot."

# Update selection to match the author's final cursor position
$ws.Range("B18").Select()
